$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 - copy formatting (bold, border, centered)
# from the neighboring "sum" header in G1 so it reuses the same cell style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New "Save" value for the single data row.
$ws.Range("H2").Value = 0
